$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

$holeIds = @(
  "LBU_05_09","LBU_05_12","LBU_87_2","LBU_05_24","LBU_05_14","LBU_05_15",
  "LBU_05_20","LBU_07_01","LBU_05_13","LBU_05_27","LBU_05_16","LBU_87_6",
  "LBU_05_25","LBU_05_29","MHZ_12_01","LBU_05_11","LBU_05_08","LBU_05_17",
  "LBU_07_02","LBU_98_4","MHZ_08_04","LBU_96_3","MHZ_08_03","MHZ_12_04",
  "LBU_02_4","MHZ_08_01","LBU_98_5","LBU_05_19","LBU_01_1","LBU_05_01",
  "LBU_96_1","LBU_87_3","LBU_87_1","MHZ_08_05","LBU_05_28","LBU_96_4",
  "LBU_87_4","LBU_05_07","MHZ_08_02","LBU_05_06","LBU_05_21","LBU_87_5",
  "LBU_02_3","LBU_05_03","MHZ_12_02","LBU_98_1"
)

# Header cell for the new index column
$ws.Range("A1").Value = "hole_id"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}
